$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks so we can rebuild them cleanly in the new row order
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-12-25 01:24:46'
$ws.Range("B2").Value = '製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5460562'
$ws.Range("G2").Value = 435
$ws.Range("H2").Value = '🔥AI,Ai ◆ツール,開発'

# Row 3
$ws.Range("A3").Value = '2025-12-25 01:24:46'
$ws.Range("B3").Value = '既存の情報検索システム(PHP)にAI文書作成システム(既存システムへの機能追加)の開発者募集します'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5460357'
$ws.Range("G3").Value = 388
$ws.Range("H3").Value = '🔥AI,Ai ◆開発 ○PHP'

# Row 4
$ws.Range("A4").Value = '2025-12-25 01:24:46'
$ws.Range("B4").Value = '産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5450864'
$ws.Range("G4").Value = 383
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

# Row 5
$ws.Range("A5").Value = '2025-12-25 01:24:46'
$ws.Range("B5").Value = '【急募】自社AIプロダクト開発|バックエンドエンジニア'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5460544'
$ws.Range("G5").Value = 375
$ws.Range("H5").Value = '🔥AI,Ai ◆開発'

# Row 6
$ws.Range("A6").Value = '2025-12-25 01:24:46'
$ws.Range("B6").Value = '【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5460294'
$ws.Range("G6").Value = 375
$ws.Range("H6").Value = '🔥AI,Ai ◆開発'

# Row 7
$ws.Range("A7").Value = '2025-12-25 01:24:46'
$ws.Range("B7").Value = '【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5460267'
$ws.Range("G7").Value = 375
$ws.Range("H7").Value = '🔥AI,Ai ◆開発'

# Row 8
$ws.Range("A8").Value = '2025-12-25 01:24:46'
$ws.Range("B8").Value = '【急募】AI活用でPDFタイトル修正のフリーランス募集!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5459721'
$ws.Range("G8").Value = 310
$ws.Range("H8").Value = '🔥AI,Ai'

# Row 9
$ws.Range("A9").Value = '2025-12-25 01:24:46'
$ws.Range("B9").Value = '施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5460563'
$ws.Range("G9").Value = 220
$ws.Range("H9").Value = '◆開発,システム開発 ◇管理'

# Row 10
$ws.Range("A10").Value = '2025-12-25 01:24:46'
$ws.Range("B10").Value = '【急募】宿泊業向けSaaSの予約者取得システム開発'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5460405'
$ws.Range("G10").Value = 118
$ws.Range("H10").Value = '◆開発,システム開発'

# Row 11
$ws.Range("A11").Value = '2025-12-25 01:24:46'
$ws.Range("B11").Value = '【急募】野球スコアボードシステム開発のフリーランス募集'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5459984'
$ws.Range("G11").Value = 118
$ws.Range("H11").Value = '◆開発,システム開発'

# Row 12
$ws.Range("A12").Value = '2025-12-25 01:24:46'
$ws.Range("B12").Value = '初回 【AWSクラウドリフト】業務アプリ移行支援エンジニア募集(Java / .NET)'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5459847'
$ws.Range("G12").Value = 103
$ws.Range("H12").Value = '★Java ◇アプリ'

# Row 13
$ws.Range("A13").Value = '2025-12-25 01:24:46'
$ws.Range("B13").Value = '現品票管理・納品書・請求書のシステムづくり'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5459942'
$ws.Range("G13").Value = 53
$ws.Range("H13").Value = '◇管理'

# Row 14
$ws.Range("A14").Value = '2025-12-25 01:24:46'
$ws.Range("B14").Value = '【急募】WEBサイト研修講師を探しています!'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5460484'
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = '◇サイト'

# Row 15
$ws.Range("A15").Value = '2025-12-25 01:24:46'
$ws.Range("B15").Value = '急募 限定公開 限定公開の仕事'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5460299'
$ws.Range("G15").Value = 18
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = '2025-12-25 01:24:46'
$ws.Range("B16").Value = '【電卓設計】ハードウェアとソフトウェアの専門家を募集!'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5459773'
$ws.Range("G16").Value = 18
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = '2025-12-25 01:24:46'
$ws.Range("B17").Value = '【電卓設計】ハードウェアとソフトウェアの専門家を募集!'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5459232'
$ws.Range("G17").Value = 18
$ws.Range("H17").ClearContents()

# Row 18
$ws.Range("A18").Value = '2025-12-25 01:24:46'
$ws.Range("B18").Value = '【急募】お名前VPSでのFTP・WPファイルアップロード改善依頼'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5459964'
$ws.Range("G18").Value = 10
$ws.Range("H18").ClearContents()

# Re-add hyperlinks for F2:F18 in row order so relationship ids come out rId1..rId17
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5460562') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5460357') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5450864') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5460544') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5460294') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5460267') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5459721') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5460563') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5460405') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5459984') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5459847') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5459942') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5460484') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5460299') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5459773') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5459232') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5459964') | Out-Null

# Restore the Hyperlink cell style (blue/underline) so style index matches s="1" instead of a duplicate style
$ws.Range("F2:F18").Style = "Hyperlink"

# Column H width: stored width 17 (ColumnWidth input includes the ~0.8333 Excel padding offset)
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668
